$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.043223783353333
$ws.Cells.Item(2, 4).Value2 = 1.060727274847082
$ws.Cells.Item(2, 5).Value2 = 1.052570876203394
$ws.Cells.Item(2, 6).Value2 = 1.066570462432922
$ws.Cells.Item(2, 9).Value2 = 1.047100223531295
$ws.Cells.Item(2, 10).Value2 = 1.048295382303557
$ws.Cells.Item(2, 11).Value2 = 1.063453267685271
$ws.Cells.Item(2, 12).Value2 = 1.055319241431795
$ws.Cells.Item(2, 13).Value2 = 1.069280655484768
$ws.Cells.Item(2, 14).Value2 = 1.049784081987589

$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.044095537052812
$ws.Cells.Item(3, 4).Value2 = 1.061377552746086
$ws.Cells.Item(3, 5).Value2 = 1.053322690009107
$ws.Cells.Item(3, 6).Value2 = 1.067348580945874
$ws.Cells.Item(3, 9).Value2 = 1.047305106140834
$ws.Cells.Item(3, 10).Value2 = 1.048814281470229
$ws.Cells.Item(3, 11).Value2 = 1.063918442096931
$ws.Cells.Item(3, 12).Value2 = 1.055884084050483
$ws.Cells.Item(3, 13).Value2 = 1.069874485409708
$ws.Cells.Item(3, 14).Value2 = 1.050303718050597

$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.044660387418126
$ws.Cells.Item(4, 4).Value2 = 1.061798862677965
$ws.Cells.Item(4, 5).Value2 = 1.053810173338565
$ws.Cells.Item(4, 6).Value2 = 1.067853034147533
$ws.Cells.Item(4, 9).Value2 = 1.047436871028215
$ws.Cells.Item(4, 10).Value2 = 1.049150146716186
$ws.Cells.Item(4, 11).Value2 = 1.064219288982384
$ws.Cells.Item(4, 12).Value2 = 1.056249916892596
$ws.Cells.Item(4, 13).Value2 = 1.07025902772874
$ws.Cells.Item(4, 14).Value2 = 1.050640060263726

$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.04489803263324
$ws.Cells.Item(5, 4).Value2 = 1.061976107930683
$ws.Cells.Item(5, 5).Value2 = 1.054015351154475
$ws.Cells.Item(5, 6).Value2 = 1.068065333739147
$ws.Cells.Item(5, 9).Value2 = 1.047492070815651
$ws.Cells.Item(5, 10).Value2 = 1.049291367783087
$ws.Cells.Item(5, 11).Value2 = 1.064345727165944
$ws.Cells.Item(5, 12).Value2 = 1.05640379341938
$ws.Cells.Item(5, 13).Value2 = 1.070420758240235
$ws.Cells.Item(5, 14).Value2 = 1.050781481880741

$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.044937944943742
$ws.Cells.Item(6, 4).Value2 = 1.062005875554129
$ws.Cells.Item(6, 5).Value2 = 1.054049815421962
$ws.Cells.Item(6, 6).Value2 = 1.068100993055287
$ws.Cells.Item(6, 9).Value2 = 1.047501327694278
$ws.Cells.Item(6, 10).Value2 = 1.049315080751489
$ws.Cells.Item(6, 11).Value2 = 1.064366954432001
$ws.Cells.Item(6, 12).Value2 = 1.056429634637823
$ws.Cells.Item(6, 13).Value2 = 1.070447917499224
$ws.Cells.Item(6, 14).Value2 = 1.050805228524278

$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.044663562131508
$ws.Cells.Item(7, 4).Value2 = 1.061801230542835
$ws.Cells.Item(7, 5).Value2 = 1.053812913994573
$ws.Cells.Item(7, 6).Value2 = 1.067855870013795
$ws.Cells.Item(7, 9).Value2 = 1.047437609374393
$ws.Cells.Item(7, 10).Value2 = 1.049152033628988
$ws.Cells.Item(7, 11).Value2 = 1.064220978606104
$ws.Cells.Item(7, 12).Value2 = 1.056251972684182
$ws.Cells.Item(7, 13).Value2 = 1.070261188510972
$ws.Cells.Item(7, 14).Value2 = 1.050641949856161

$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.043518236731487
$ws.Cells.Item(8, 4).Value2 = 1.060946927038815
$ws.Cells.Item(8, 5).Value2 = 1.052824745053279
$ws.Cells.Item(8, 6).Value2 = 1.066833231563413
$ws.Cells.Item(8, 9).Value2 = 1.047169631379358
$ws.Cells.Item(8, 10).Value2 = 1.04847072464771
$ws.Cells.Item(8, 11).Value2 = 1.063610506174207
$ws.Cells.Item(8, 12).Value2 = 1.05551006097464
$ws.Cells.Item(8, 13).Value2 = 1.06948128110568
$ws.Cells.Item(8, 14).Value2 = 1.049959673337983

$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.041505975644019
$ws.Cells.Item(9, 4).Value2 = 1.059445729963601
$ws.Cells.Item(9, 5).Value2 = 1.051091279055752
$ws.Cells.Item(9, 6).Value2 = 1.065038641294987
$ws.Cells.Item(9, 9).Value2 = 1.046691266683192
$ws.Cells.Item(9, 10).Value2 = 1.047271015490797
$ws.Cells.Item(9, 11).Value2 = 1.062533674936896
$ws.Cells.Item(9, 12).Value2 = 1.054205400473549
$ws.Cells.Item(9, 13).Value2 = 1.068109312895829
$ws.Cells.Item(9, 14).Value2 = 1.048758260456458

$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.040168563314678
$ws.Cells.Item(10, 4).Value2 = 1.058447866865425
$ws.Cells.Item(10, 5).Value2 = 1.049940999689803
$ws.Cells.Item(10, 6).Value2 = 1.063847362325732
$ws.Cells.Item(10, 9).Value2 = 1.046368265744039
$ws.Cells.Item(10, 10).Value2 = 1.046471854047092
$ws.Cells.Item(10, 11).Value2 = 1.061815133530938
$ws.Cells.Item(10, 12).Value2 = 1.053337517160908
$ws.Cells.Item(10, 13).Value2 = 1.067196329555817
$ws.Cells.Item(10, 14).Value2 = 1.047957964111839

$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.03959044039557
$ws.Cells.Item(11, 4).Value2 = 1.058016501524154
$ws.Cells.Item(11, 5).Value2 = 1.049444211924151
$ws.Cells.Item(11, 6).Value2 = 1.063332763769041
$ws.Cells.Item(11, 9).Value2 = 1.046227442105494
$ws.Cells.Item(11, 10).Value2 = 1.04612597746504
$ws.Cells.Item(11, 11).Value2 = 1.061503860303767
$ws.Cells.Item(11, 12).Value2 = 1.052962180777694
$ws.Cells.Item(11, 13).Value2 = 1.066801412379607
$ws.Cells.Item(11, 14).Value2 = 1.047611596345368

$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.03937584946521
$ws.Cells.Item(12, 4).Value2 = 1.057856382766944
$ws.Cells.Item(12, 5).Value2 = 1.049259878663453
$ws.Cells.Item(12, 6).Value2 = 1.063141806266966
$ws.Cells.Item(12, 9).Value2 = 1.046174990161497
$ws.Cells.Item(12, 10).Value2 = 1.045997529687504
$ws.Cells.Item(12, 11).Value2 = 1.06138821999681
$ws.Cells.Item(12, 12).Value2 = 1.052822835145573
$ws.Cells.Item(12, 13).Value2 = 1.066654785904265
$ws.Cells.Item(12, 14).Value2 = 1.047482966157254

$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.039421873149756
$ws.Cells.Item(13, 4).Value2 = 1.057890723795492
$ws.Cells.Item(13, 5).Value2 = 1.049299409881038
$ws.Cells.Item(13, 6).Value2 = 1.063182758790026
$ws.Cells.Item(13, 9).Value2 = 1.046186247775248
$ws.Cells.Item(13, 10).Value2 = 1.046025080965895
$ws.Cells.Item(13, 11).Value2 = 1.061413026097215
$ws.Cells.Item(13, 12).Value2 = 1.05285272201947
$ws.Cells.Item(13, 13).Value2 = 1.066686234887791
$ws.Cells.Item(13, 14).Value2 = 1.04751055656162

$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.039572699179917
$ws.Cells.Item(14, 4).Value2 = 1.058003263811819
$ws.Cells.Item(14, 5).Value2 = 1.049428970873314
$ws.Cells.Item(14, 6).Value2 = 1.063316975329312
$ws.Cells.Item(14, 9).Value2 = 1.046223109339934
$ws.Cells.Item(14, 10).Value2 = 1.046115359389952
$ws.Cells.Item(14, 11).Value2 = 1.061494301830654
$ws.Cells.Item(14, 12).Value2 = 1.052950660967086
$ws.Cells.Item(14, 13).Value2 = 1.066789290887577
$ws.Cells.Item(14, 14).Value2 = 1.047600963191395

$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.039665647955014
$ws.Cells.Item(15, 4).Value2 = 1.058072618003815
$ws.Cells.Item(15, 5).Value2 = 1.049508823683795
$ws.Cells.Item(15, 6).Value2 = 1.063399695451297
$ws.Cells.Item(15, 9).Value2 = 1.046245801938083
$ws.Cells.Item(15, 10).Value2 = 1.046170986415457
$ws.Cells.Item(15, 11).Value2 = 1.061544375943117
$ws.Cells.Item(15, 12).Value2 = 1.053011013839191
$ws.Cells.Item(15, 13).Value2 = 1.06685279553397
$ws.Cells.Item(15, 14).Value2 = 1.047656669213656

$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.040206952328395
$ws.Cells.Item(16, 4).Value2 = 1.058476510415356
$ws.Cells.Item(16, 5).Value2 = 1.049973997216066
$ws.Cells.Item(16, 6).Value2 = 1.063881540692583
$ws.Cells.Item(16, 9).Value2 = 1.04637759154049
$ws.Cells.Item(16, 10).Value2 = 1.046494812337126
$ws.Cells.Item(16, 11).Value2 = 1.061835788884608
$ws.Cells.Item(16, 12).Value2 = 1.053362436877838
$ws.Cells.Item(16, 13).Value2 = 1.067222547700999
$ws.Cells.Item(16, 14).Value2 = 1.047980955005278

$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.040546763227963
$ws.Cells.Item(17, 4).Value2 = 1.058730054531083
$ws.Cells.Item(17, 5).Value2 = 1.050266135211875
$ws.Cells.Item(17, 6).Value2 = 1.064184121068992
$ws.Cells.Item(17, 9).Value2 = 1.046460002637761
$ws.Cells.Item(17, 10).Value2 = 1.046697985039093
$ws.Cells.Item(17, 11).Value2 = 1.062018548065127
$ws.Cells.Item(17, 12).Value2 = 1.05358300021935
$ws.Cells.Item(17, 13).Value2 = 1.067454594573212
$ws.Cells.Item(17, 14).Value2 = 1.048184416235785

$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.04074506414942
$ws.Cells.Item(18, 4).Value2 = 1.058878011472336
$ws.Cells.Item(18, 5).Value2 = 1.050436658776901
$ws.Cells.Item(18, 6).Value2 = 1.064360729921629
$ws.Cells.Item(18, 9).Value2 = 1.046507978756012
$ws.Cells.Item(18, 10).Value2 = 1.046816508128769
$ws.Cells.Item(18, 11).Value2 = 1.062125134791002
$ws.Cells.Item(18, 12).Value2 = 1.053711695640868
$ws.Cells.Item(18, 13).Value2 = 1.067589983097221
$ws.Cells.Item(18, 14).Value2 = 1.048303107641844

$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.040812695699196
$ws.Cells.Item(19, 4).Value2 = 1.058928472600549
$ws.Cells.Item(19, 5).Value2 = 1.050494823950662
$ws.Cells.Item(19, 6).Value2 = 1.064420969070892
$ws.Cells.Item(19, 9).Value2 = 1.046524321594973
$ws.Cells.Item(19, 10).Value2 = 1.046856924094941
$ws.Cells.Item(19, 11).Value2 = 1.062161475756477
$ws.Cells.Item(19, 12).Value2 = 1.053755584966295
$ws.Cells.Item(19, 13).Value2 = 1.067636153736296
$ws.Cells.Item(19, 14).Value2 = 1.048343581003323

$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.040510294899265
$ws.Cells.Item(20, 4).Value2 = 1.05870284450113
$ws.Cells.Item(20, 5).Value2 = 1.050234778687135
$ws.Cells.Item(20, 6).Value2 = 1.064151644751365
$ws.Cells.Item(20, 9).Value2 = 1.046451170308469
$ws.Cells.Item(20, 10).Value2 = 1.046676184877976
$ws.Cells.Item(20, 11).Value2 = 1.061998941132025
$ws.Cells.Item(20, 12).Value2 = 1.053559331240884
$ws.Cells.Item(20, 13).Value2 = 1.067429694041321
$ws.Cells.Item(20, 14).Value2 = 1.048162585115939

$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.039528280528893
$ws.Cells.Item(21, 4).Value2 = 1.057970120542418
$ws.Cells.Item(21, 5).Value2 = 1.04939081296899
$ws.Cells.Item(21, 6).Value2 = 1.063277446716459
$ws.Cells.Item(21, 9).Value2 = 1.046212258490568
$ws.Cells.Item(21, 10).Value2 = 1.046088773910864
$ws.Cells.Item(21, 11).Value2 = 1.061470368681343
$ws.Cells.Item(21, 12).Value2 = 1.052921818415907
$ws.Cells.Item(21, 13).Value2 = 1.066758941702256
$ws.Cells.Item(21, 14).Value2 = 1.047574339957877

$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.038911715821351
$ws.Cells.Item(22, 4).Value2 = 1.057510062629597
$ws.Cells.Item(22, 5).Value2 = 1.048861312209377
$ws.Cells.Item(22, 6).Value2 = 1.062728888474736
$ws.Cells.Item(22, 9).Value2 = 1.046061213696234
$ws.Cells.Item(22, 10).Value2 = 1.045719597468348
$ws.Cells.Item(22, 11).Value2 = 1.061137921699077
$ws.Cells.Item(22, 12).Value2 = 1.052521400388441
$ws.Cells.Item(22, 13).Value2 = 1.066337580412876
$ws.Cells.Item(22, 14).Value2 = 1.047204639242467

$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.039238485746059
$ws.Cells.Item(23, 4).Value2 = 1.057753887211866
$ws.Cells.Item(23, 5).Value2 = 1.049141902455045
$ws.Cells.Item(23, 6).Value2 = 1.063019586069113
$ws.Cells.Item(23, 9).Value2 = 1.046141363953908
$ws.Cells.Item(23, 10).Value2 = 1.045915290027392
$ws.Cells.Item(23, 11).Value2 = 1.061314168307416
$ws.Cells.Item(23, 12).Value2 = 1.052733630053195
$ws.Cells.Item(23, 13).Value2 = 1.066560916633999
$ws.Cells.Item(23, 14).Value2 = 1.047400609707392

$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.040526773073983
$ws.Cells.Item(24, 4).Value2 = 1.058715139329646
$ws.Cells.Item(24, 5).Value2 = 1.050248946966767
$ws.Cells.Item(24, 6).Value2 = 1.064166319034452
$ws.Cells.Item(24, 9).Value2 = 1.046455161544974
$ws.Cells.Item(24, 10).Value2 = 1.046686035383472
$ws.Cells.Item(24, 11).Value2 = 1.062007800703372
$ws.Cells.Item(24, 12).Value2 = 1.05357002609595
$ws.Cells.Item(24, 13).Value2 = 1.067440945397987
$ws.Cells.Item(24, 14).Value2 = 1.048172449610283

$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.042025478599619
$ws.Cells.Item(25, 4).Value2 = 1.059833316466239
$ws.Cells.Item(25, 5).Value2 = 1.051538484136583
$ws.Cells.Item(25, 6).Value2 = 1.065501693166069
$ws.Cells.Item(25, 9).Value2 = 1.046815659886417
$ws.Cells.Item(25, 10).Value2 = 1.047581061079967
$ws.Cells.Item(25, 11).Value2 = 1.062812183191791
$ws.Cells.Item(25, 12).Value2 = 1.054542359642004
$ws.Cells.Item(25, 13).Value2 = 1.068463713683119
$ws.Cells.Item(25, 14).Value2 = 1.049068746345928
